# Update countries & provincias Spain
# Apply data refresh to the "Pais" worksheet:
#   1) Estados Unidos (row 4): refresh "Casos activos" / "Recuperados"
#   2) Bolivia's case counts increased enough to move it up the ranking,
#      ahead of Eslovenia/Eslovaquia/Lituania/Costa de Marfil (rows 84-88)
#   3) San Cristobal y Nieves overtakes Burundi in the ranking (rows 198-199)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Estados Unidos row (row 4): update Casos activos / Recuperados ---
$ws.Cells.Item(4, 4).Value = 173725   # D4 Casos activos
$ws.Cells.Item(4, 5).Value = 919605   # E4 Recuperados

# --- 2) Re-rank rows 84-88 ---
# Before:
#   84 Eslovenia        1439 0 239 1106 21 0 94
#   85 Eslovaquia        1407 0 608  775  5 0 24
#   86 Lituania          1406 0 632  728 17 0 46
#   87 Costa de Marfil   1362 0 622  725  0 0 15
#   88 Bolivia           1229 0 134 1029  3 0 66
# After (Bolivia's updated totals push it to rank 84; the rest shift down one row):
#   84 Bolivia           1470 241 159 1240 3 5 71
#   85 Eslovenia         1439   0 239 1106 21 0 94
#   86 Eslovaquia        1407   0 608  775  5 0 24
#   87 Lituania          1406   0 632  728 17 0 46
#   88 Costa de Marfil   1362   0 622  725  0 0 15

$ws.Cells.Item(84, 1).Value = "Bolivia"
$ws.Cells.Item(84, 2).Value = 1470
$ws.Cells.Item(84, 3).Value = 241
$ws.Cells.Item(84, 4).Value = 159
$ws.Cells.Item(84, 5).Value = 1240
$ws.Cells.Item(84, 6).Value = 3
$ws.Cells.Item(84, 7).Value = 5
$ws.Cells.Item(84, 8).Value = 71

$ws.Cells.Item(85, 1).Value = "Eslovenia"
$ws.Cells.Item(85, 2).Value = 1439
$ws.Cells.Item(85, 3).Value = 0
$ws.Cells.Item(85, 4).Value = 239
$ws.Cells.Item(85, 5).Value = 1106
$ws.Cells.Item(85, 6).Value = 21
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 94

$ws.Cells.Item(86, 1).Value = "Eslovaquia"
$ws.Cells.Item(86, 2).Value = 1407
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 4).Value = 608
$ws.Cells.Item(86, 5).Value = 775
$ws.Cells.Item(86, 6).Value = 5
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 24

$ws.Cells.Item(87, 1).Value = "Lituania"
$ws.Cells.Item(87, 2).Value = 1406
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = 632
$ws.Cells.Item(87, 5).Value = 728
$ws.Cells.Item(87, 6).Value = 17
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 46

$ws.Cells.Item(88, 1).Value = "Costa de Marfil"
$ws.Cells.Item(88, 2).Value = 1362
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 622
$ws.Cells.Item(88, 5).Value = 725
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 15

# Row 89 (Senegal) is unaffected by the re-rank and keeps its values.

# --- 3) Swap Burundi / San Cristobal y Nieves (rows 198-199) ---
# Before: 198 Burundi (7 activos / 1 muerte), 199 San Cristobal y Nieves (8 activos / 0 muertes)
# After:  198 San Cristobal y Nieves (8 activos / 0 muertes), 199 Burundi (7 activos / 1 muerte)
$ws.Cells.Item(198, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(198, 4).Value = 8
$ws.Cells.Item(198, 8).Value = 0

$ws.Cells.Item(199, 1).Value = "Burundi"
$ws.Cells.Item(199, 4).Value = 7
$ws.Cells.Item(199, 8).Value = 1
